# "Falta descargar archivos en splash"
# Adds a new "iniciales" worksheet (between "signos1" and "Datos Signos")
# that lists the files ("archivos"/"imágenes") that must be downloaded
# before showing the splash screen.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new worksheet right after "signos1" -------------------
$afterSheet = $wb.Worksheets.Item("signos1")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "iniciales"

# --- 2. Fill in the data --------------------------------------------------
# Column C (file names) is populated first, then column B (group labels),
# then column A (index within the group) and finally column D (formula),
# matching the order the original sheet's shared strings were recorded in.

$fileNames = @(
    "compatibilidad.json",
    "signos.json",
    "acuario.png",
    "aries.png",
    "cancer.png",
    "capricornio.png",
    "escorpio.png",
    "geminis.png",
    "leo.png",
    "libra.png",
    "piscis.png",
    "sagitario.png",
    "taruro.png",
    "virgo.png"
)

for ($i = 0; $i -lt $fileNames.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $fileNames[$i]
}

$ws.Cells.Item(1, 2).Value = "archivos"
$ws.Cells.Item(2, 2).Value = "archivos"
for ($row = 3; $row -le 14; $row++) {
    $ws.Cells.Item($row, 2).Value = "imágenes"
}

$ws.Cells.Item(1, 1).Value = 0
$ws.Cells.Item(2, 1).Value = 1
for ($row = 3; $row -le 14; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 3
}

for ($row = 1; $row -le 14; $row++) {
    $ws.Cells.Item($row, 4).Formula = '=B' + $row + '&"["&A' + $row + '&"]   archivo="&C' + $row
}

# --- 3. Column widths / row height ---------------------------------------
$ws.Columns.Item(3).ColumnWidth = 16.1      # -> width 17, bestFit-like
$ws.Columns.Item(4).ColumnWidth = 61.3      # -> width ~62.1640625
$ws.Rows.Item(14).RowHeight = 17

# --- 4. Activate the sheet, set zoom & selection --------------------------
$ws.Activate() | Out-Null
$ws.Range("D1:D14").Select() | Out-Null
$excel.ActiveWindow.Zoom = 140

Write-Host "iniciales sheet created"
